$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column I ("ISI") - reuse header H's style (bold, left-aligned)
$ws.Range("I1").Value = "ISI"
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill data rows 2-52 of column I with the ISI values, grouped by protocol block
for ($r = 2; $r -le 20; $r++) {
    $ws.Cells.Item($r, 9).Value = 4
}
for ($r = 21; $r -le 38; $r++) {
    $ws.Cells.Item($r, 9).Value = 5
}
for ($r = 39; $r -le 52; $r++) {
    $ws.Cells.Item($r, 9).Value = 2
}

# Restore the user's last selection
$ws.Range("M8").Select()
